# edit.ps1
# Refresh the "cryptos" price list: update the Price (D) and Volume(1h) (E)
# columns for rows 2-51 of Sheet1 to the latest scraped values.
#
# Columns D and E hold plain text (inline string) cells in the original
# workbook, not numbers. Several of the new Price strings (e.g. "323.40",
# "1.001") are valid decimal numbers, so a bare .Value assignment would be
# auto-converted to a numeric cell by Excel. To avoid that, each D-column
# write briefly forces the cell to Text format, assigns the value, then
# restores the original General/Normal formatting so the cell ends up
# with no stray style, exactly like the source cells.
#
# NOTE: this is written as flat, top-level statements (no helper
# functions/cmdlets) because passing range-address / value strings through
# PowerShell function parameters into COM property assignments in this
# runtime does not reliably apply the write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.379.75"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.89"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("E4").Value = "  -1.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.40"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4533"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3871"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.80"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07930"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.72%  "
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.40"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.874.36"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.921"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.129"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.00"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06546"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.03"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.49%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.531"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.380.61"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.84"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.085.32"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.87"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.83"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.077"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.443"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.76%  "
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.487"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09290"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9381"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.597"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.261"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02238"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.224"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05993"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.199"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.64%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5915"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1894"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.13"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.276"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5614"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.05"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.375"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.925"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06767"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.42"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.51%  "
